# EIA Table 1.6.B update: October 2016 YTD -> November 2016 YTD
# (commit: "2017-01-31 update: energy.gov - chunk 7")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Title / header text: "October" -> "November"
# ---------------------------------------------------------------
$ws.Range("A2").Value = "by State, by Sector, Year-to-Date through November 2016 and 2015 (Thousand Megawatthours)"

# Column-pair headers in row 6 (repeated across the All/Electric/Commercial/
# Industrial/Electric Utilities/Independent Power Producers column groups)
$ws.Range("B6").Value = "November 2016 YTD"
$ws.Range("E6").Value = "November 2016 YTD"
$ws.Range("G6").Value = "November 2016 YTD"
$ws.Range("I6").Value = "November 2016 YTD"
$ws.Range("K6").Value = "November 2016 YTD"

$ws.Range("C6").Value = "November 2015 YTD"
$ws.Range("F6").Value = "November 2015 YTD"
$ws.Range("H6").Value = "November 2015 YTD"
$ws.Range("J6").Value = "November 2015 YTD"
$ws.Range("L6").Value = "November 2015 YTD"

# ---------------------------------------------------------------
# 2) Updated data values (YTD figures revised for the November release)
# ---------------------------------------------------------------

function Set-RowValues {
    param($ws, $row, $values)
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

Set-RowValues $ws 14 @{ B=241; C=211; D=0.14000000000000001; K=241; L=211 }
Set-RowValues $ws 15 @{ B=68; C=65; D=0.033; K=68; L=65 }
Set-RowValues $ws 17 @{ B=173; C=146; D=0.189; K=173; L=146 }
Set-RowValues $ws 18 @{ B=2160; C=3000; D=-0.28; E=1094; F=1659; G=875; H=1091; K=191; L=250 }
Set-RowValues $ws 20 @{ C=1060; D=-0.531; F=1060 }
Set-RowValues $ws 21 @{ B=664; C=703; D=-0.055; E=519; F=542; H=25; K=142; L=137 }
Set-RowValues $ws 22 @{ B=885; C=1075; D=-0.176; G=872; H=1066; L=9 }
Set-RowValues $ws 23 @{ B=113; C=162; D=-0.301; E=78; F=58; K=35; L=104 }
Set-RowValues $ws 24 @{ B=57; C=43; D=0.319; J=8; K=52; L=35 }
Set-RowValues $ws 25 @{ B=57; C=43; D=0.319; J=8; K=52; L=35 }
Set-RowValues $ws 32 @{ B=2049; C=1570; D=0.305; E=1952; F=1464; K=96; L=105 }
Set-RowValues $ws 35 @{ B=1952; C=1464; D=0.33300000000000002; E=1952; F=1464 }
Set-RowValues $ws 36 @{ B=96; C=105; D=-0.085; K=96; L=105 }
Set-RowValues $ws 42 @{ B=1051; C=906; D=0.16; E=1051; F=906 }
Set-RowValues $ws 44 @{ B=1051; C=906; D=0.16; E=1051; F=906 }
Set-RowValues $ws 47 @{ B=4407; C=3955; D=0.114; E=4116; F=3644; K=291; L=311 }
Set-RowValues $ws 49 @{ B=4296; C=3829; D=0.122; E=4116; F=3644; K=180; L=185 }
Set-RowValues $ws 51 @{ B=111; C=126; D=-0.119; K=111; L=126 }
Set-RowValues $ws 52 @{ B=401; C=442; D=-0.093; G=401; H=442 }
Set-RowValues $ws 56 @{ B=401; C=442; D=-0.093; G=401; H=442 }
Set-RowValues $ws 68 @{ B=10366; C=10128; D=0.023; E=8214; F=7674; G=1276; H=1533; J=8; K=871; L=912 }
